# Fix inconsistent/missing contributor names (per commit message):
#  - "Georg P Krog"  -> "Georg P. Krog"   (add missing period after middle initial)
#  - "Georg P Krogg" -> "Georg P. Krog"   (fix typo'd surname + add period)
# Applied across every cell (column N = Contributors) in every sheet that
# holds one of the affected strings, so the shared-string table stays
# correctly de-duplicated.

$wb = $excel.ActiveWorkbook

$replacements = @(
    @{ Old = "Georg P Krog, Harshvardhan J. Pandit, Paul Ryan"; New = "Georg P. Krog, Harshvardhan J. Pandit, Paul Ryan" },
    @{ Old = "David Hickey, Georg P Krogg"; New = "David Hickey, Georg P. Krog" },
    @{ Old = "Georg P Krog"; New = "Georg P. Krog" },
    @{ Old = "Georg P Krog, Harshvardhan J. Pandit, Paul Ryan, Julian Flake"; New = "Georg P. Krog, Harshvardhan J. Pandit, Paul Ryan, Julian Flake" }
)

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($row in $used.Rows) {
        $cell = $ws.Cells.Item($row.Row, 14)  # column N = Contributors
        $val = $cell.Value()
        if ($val -eq $null) { continue }
        foreach ($rep in $replacements) {
            if ($val -eq $rep.Old) {
                $cell.Value = $rep.New
                break
            }
        }
    }
}

# The LegalBasis sheet's "ContractPerformance / EnterIntoContract /
# LegitimateInterestOfController / LegitimateInterestOfThirdParty /
# OfficialAuthorityOfController / VitalInterestOfDataSubject /
# VitalInterestOfNaturalPerson" rows shared a (now corrected) Contributors
# cell style that was out of step with the other Contributors cells
# (e.g. N10, N16). Normalize their formatting to match by pasting the
# format from an already-consistent cell (N10), leaving cell values intact.
$wsLegalBasis = $wb.Worksheets.Item("LegalBasis")
$formatSource = $wsLegalBasis.Range("N10")
$formatSource.Copy()
foreach ($addr in @("N9", "N11", "N14", "N15", "N17", "N20", "N21")) {
    $wsLegalBasis.Range($addr).PasteSpecial(-4122)
}

Write-Output "contributors fixed"
